$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 929 (shifts old 929-990 down to 934-995)
$ws.Rows.Item(929).Resize(5).Insert()

# Row 929
$ws.Cells.Item(929,1).Value2 = 6
$ws.Cells.Item(929,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(929,3).Value = "Metropolitana"
$ws.Cells.Item(929,4).Value2 = 44585
$ws.Cells.Item(929,5).Value2 = 13
$ws.Cells.Item(929,6).Value2 = 100112027
$ws.Cells.Item(929,7).Value = "Melón"
$ws.Cells.Item(929,8).Value = "Tuna"
$ws.Cells.Item(929,9).Value = "Extra"
$ws.Cells.Item(929,10).Value2 = 4900
$ws.Cells.Item(929,11).Value2 = 700
$ws.Cells.Item(929,12).Value2 = 800
$ws.Cells.Item(929,13).Value2 = 745
$ws.Cells.Item(929,14).Value = "`$/unidad"
$ws.Cells.Item(929,15).Value = "Región de O'Higgins"
$ws.Cells.Item(929,16).Value2 = 745
$ws.Cells.Item(929,17).Value2 = 1
$ws.Cells.Item(929,18).Value = "Hortaliza"

# Row 930
$ws.Cells.Item(930,1).Value2 = 6
$ws.Cells.Item(930,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(930,3).Value = "Metropolitana"
$ws.Cells.Item(930,4).Value2 = 44585
$ws.Cells.Item(930,5).Value2 = 13
$ws.Cells.Item(930,6).Value2 = 100112027
$ws.Cells.Item(930,7).Value = "Melón"
$ws.Cells.Item(930,8).Value = "Tuna"
$ws.Cells.Item(930,9).Value = "Primera"
$ws.Cells.Item(930,10).Value2 = 6700
$ws.Cells.Item(930,11).Value2 = 500
$ws.Cells.Item(930,12).Value2 = 600
$ws.Cells.Item(930,13).Value2 = 543
$ws.Cells.Item(930,14).Value = "`$/unidad"
$ws.Cells.Item(930,15).Value = "Región de O'Higgins"
$ws.Cells.Item(930,16).Value2 = 543
$ws.Cells.Item(930,17).Value2 = 1
$ws.Cells.Item(930,18).Value = "Hortaliza"

# Row 931
$ws.Cells.Item(931,1).Value2 = 6
$ws.Cells.Item(931,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(931,3).Value = "Metropolitana"
$ws.Cells.Item(931,4).Value2 = 44585
$ws.Cells.Item(931,5).Value2 = 13
$ws.Cells.Item(931,6).Value2 = 100112027
$ws.Cells.Item(931,7).Value = "Melón"
$ws.Cells.Item(931,8).Value = "Tuna"
$ws.Cells.Item(931,9).Value = "Segunda"
$ws.Cells.Item(931,10).Value2 = 6200
$ws.Cells.Item(931,11).Value2 = 350
$ws.Cells.Item(931,12).Value2 = 450
$ws.Cells.Item(931,13).Value2 = 390
$ws.Cells.Item(931,14).Value = "`$/unidad"
$ws.Cells.Item(931,15).Value = "Región de O'Higgins"
$ws.Cells.Item(931,16).Value2 = 390
$ws.Cells.Item(931,17).Value2 = 1
$ws.Cells.Item(931,18).Value = "Hortaliza"

# Row 932
$ws.Cells.Item(932,1).Value2 = 6
$ws.Cells.Item(932,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(932,3).Value = "Metropolitana"
$ws.Cells.Item(932,4).Value2 = 44585
$ws.Cells.Item(932,5).Value2 = 13
$ws.Cells.Item(932,6).Value2 = 100112027
$ws.Cells.Item(932,7).Value = "Melón"
$ws.Cells.Item(932,8).Value = "Tuna"
$ws.Cells.Item(932,9).Value = "Super"
$ws.Cells.Item(932,10).Value2 = 3800
$ws.Cells.Item(932,11).Value2 = 900
$ws.Cells.Item(932,12).Value2 = 1000
$ws.Cells.Item(932,13).Value2 = 945
$ws.Cells.Item(932,14).Value = "`$/unidad"
$ws.Cells.Item(932,15).Value = "Región de O'Higgins"
$ws.Cells.Item(932,16).Value2 = 945
$ws.Cells.Item(932,17).Value2 = 1
$ws.Cells.Item(932,18).Value = "Hortaliza"

# Row 933
$ws.Cells.Item(933,1).Value2 = 6
$ws.Cells.Item(933,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(933,3).Value = "Metropolitana"
$ws.Cells.Item(933,4).Value2 = 44585
$ws.Cells.Item(933,5).Value2 = 13
$ws.Cells.Item(933,6).Value2 = 100112027
$ws.Cells.Item(933,7).Value = "Melón"
$ws.Cells.Item(933,8).Value = "Tuna"
$ws.Cells.Item(933,9).Value = "Tercera"
$ws.Cells.Item(933,10).Value2 = 4900
$ws.Cells.Item(933,11).Value2 = 200
$ws.Cells.Item(933,12).Value2 = 300
$ws.Cells.Item(933,13).Value2 = 245
$ws.Cells.Item(933,14).Value = "`$/unidad"
$ws.Cells.Item(933,15).Value = "Región de O'Higgins"
$ws.Cells.Item(933,16).Value2 = 245
$ws.Cells.Item(933,17).Value2 = 1
$ws.Cells.Item(933,18).Value = "Hortaliza"
